# Update the "selected_features" sheet: rename the *_corr column headers,
# and rewrite the feature-ranking lists for columns A-F (rows 2-15) to
# reflect the new ridge-regression / random-forest results, extending the
# used range from A1:F11 to A1:F15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("D1").Value = 'rfe_corr'
$ws.Range("E1").Value = 'mutual information_corr'
$ws.Range("F1").Value = 'random forest_corr'

# Row 3
$ws.Range("C3").Value = 'Pre-Test-Cell-Open-Circuit-Voltage-V'
$ws.Range("D3").Value = 'Pre-Test-Cell-Open-Circuit-Voltage-V'
$ws.Range("F3").Value = 'Pre-Test-Cell-Open-Circuit-Voltage-V'

# Row 4
$ws.Range("C4").Value = 'Bottom-Vent-Yes-No'
$ws.Range("D4").Value = 'Pre-Test-Cell-Mass-g'
$ws.Range("E4").Value = 'Pre-Test-Cell-Open-Circuit-Voltage-V'
$ws.Range("F4").Value = 'Bottom-Vent-Yes-No'

# Row 5
$ws.Range("A5").Value = 'State-of-Charge'
$ws.Range("B5").Value = 'Pre-Test-Cell-Open-Circuit-Voltage-V'
$ws.Range("C5").Value = 'State-of-Charge'
$ws.Range("D5").Value = 'Trigger-Mechanism_Nail'
$ws.Range("E5").Value = 'Cell-Nominal-Voltage-V'
$ws.Range("F5").Value = 'Cell-Capacity-Ah'

# Row 6
$ws.Range("A6").Value = 'Trigger-Mechanism_Nail'
$ws.Range("B6").Value = 'State-of-Charge'
$ws.Range("C6").Value = 'Cell-Energy-Wh'
$ws.Range("D6").Value = 'Chemistry_NCA/Graphite'
$ws.Range("E6").Value = 'Bottom-Vent-Yes-No'
$ws.Range("F6").Value = 'Trigger-Mechanism_Nail'

# Row 7
$ws.Range("A7").Value = 'Chemistry_NCA/Graphite'
$ws.Range("B7").Value = 'Cell-Nominal-Voltage-V'
$ws.Range("C7").Value = 'Cell-Capacity-Ah'
$ws.Range("D7").Value = 'Chemistry_NMC/Graphite'
$ws.Range("E7").Value = 'Chemistry_NMC/Graphite'
$ws.Range("F7").Value = 'Cell-Nominal-Voltage-V'

# Row 8
$ws.Range("A8").Value = 'Pressure-Assisted-Seal-Configuration-Negative'
$ws.Range("B8").Value = 'Bottom-Vent-Yes-No'
$ws.Range("C8").Value = 'Trigger-Mechanism_Nail'
$ws.Range("D8").Value = 'Pressure-Assisted-Seal-Configuration-Negative'
$ws.Range("E8").Value = 'Trigger-Mechanism_Nail'

# Row 9
$ws.Range("A9").Value = 'Bottom-Vent-Yes-No'
$ws.Range("B9").Value = 'Chemistry_NMC/Graphite'
$ws.Range("C9").Value = 'Cell-Nominal-Voltage-V'
$ws.Range("D9").Value = 'Bottom-Vent-Yes-No'
$ws.Range("E9").Value = 'Chemistry_NCA/Graphite'
$ws.Range("F9").Value = 'Trigger-Mechanism_Heater (Non-ISC)'

# Row 10
$ws.Range("B10").Value = 'Trigger-Mechanism_Nail'
$ws.Range("C10").Value = 'Pressure-Assisted-Seal-Configuration-Positive'
$ws.Range("E10").Value = 'Trigger-Mechanism_Heater (Non-ISC)'
$ws.Range("F10").Value = 'Pressure-Assisted-Seal-Configuration-Negative'

# Row 11
$ws.Range("B11").Value = 'Chemistry_NCA/Graphite'
$ws.Range("C11").Value = 'Trigger-Mechanism_Heater (Non-ISC)'
$ws.Range("E11").Value = 'Trigger-Mechanism_Heater (ISC)'

# Row 12
$ws.Range("B12").Value = 'Trigger-Mechanism_Heater (Non-ISC)'
$ws.Range("C12").Value = 'Pressure-Assisted-Seal-Configuration-Negative'
$ws.Range("E12").Value = 'Pressure-Assisted-Seal-Configuration-Positive'

# Row 13
$ws.Range("B13").Value = 'Trigger-Mechanism_Heater (ISC)'
$ws.Range("C13").Value = 'Chemistry_NCA/Graphite'
$ws.Range("E13").Value = 'Pressure-Assisted-Seal-Configuration-Negative'

# Row 14
$ws.Range("B14").Value = 'Pressure-Assisted-Seal-Configuration-Positive'

# Row 15
$ws.Range("B15").Value = 'Pressure-Assisted-Seal-Configuration-Negative'
